$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data table (Cohort, Component, Direction, Mean, Standard_Deviation, Participants).
# Rows 2-13 are the recalculated Mean/Standard_Deviation values for the existing
# 100s/200s/400s cohorts, and rows 14-17 are the newly added 500s cohort data.
$rows = @(
  @("100s", "Rambling", "X", 0.03798002009230238, 0.3675995650866852, 5),
  @("100s", "Rambling", "Y", -0.4898621027432221, 2.468738120537818, 5),
  @("100s", "Trembling", "X", -0.0004364931461946009, 0.01892820660228399, 5),
  @("100s", "Trembling", "Y", -0.1625615005501911, 1.423318893505352, 5),
  @("200s", "Rambling", "X", 0.09329098862128271, 0.09138359560488801, 1),
  @("200s", "Rambling", "Y", 0.9256992358248586, 3.986278789203952, 1),
  @("200s", "Trembling", "X", -0.001912460049854115, 0.01862014615541836, 1),
  @("200s", "Trembling", "Y", 0.2321302713179986, 3.911545601050795, 1),
  @("400s", "Rambling", "X", -0.1259622366904592, 0.4172340946629315, 15),
  @("400s", "Rambling", "Y", 0.1877171198070554, 4.459132562394634, 15),
  @("400s", "Trembling", "X", -0.0002185402936161674, 0.03108289684252993, 15),
  @("400s", "Trembling", "Y", -0.1137896234410918, 3.761602244731122, 15),
  @("500s", "Rambling", "X", 0.1611015679084914, 0.05824351391409256, 1),
  @("500s", "Rambling", "Y", 2.281324606821694, 1.31477609252953, 1),
  @("500s", "Trembling", "X", 0.007713486145562705, 0.02003382936956402, 1),
  @("500s", "Trembling", "Y", 0.2847299675026307, 0.9918943620372968, 1)
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r++
}
